# Update: pi 08. 01. 2021
# Inserts a new block of 71 district rows for date 2021-01-06 (serial 44202)
# immediately before the existing 2021-01-07 (serial 44203) block, shifting
# the existing data down by 71 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 7633
$newDate = 44202

# District name / daily positive-test count pairs for 2021-01-06.
$districts = @(
    @('Bánovce nad Bebravou', 9),
    @('Banská Bystrica', 135),
    @('Banská Štiavnica', 4),
    @('Bardejov', 33),
    @('Bratislava', 172),
    @('Brezno', 19),
    @('Bytča', 6),
    @('Čadca', 54),
    @('Detva', 8),
    @('Dolný Kubín', 6),
    @('Dunajská Streda', 79),
    @('Galanta', 44),
    @('Gelnica', 7),
    @('Hlohovec', 10),
    @('Humenné', 5),
    @('Ilava', 74),
    @('Kežmarok', 36),
    @('Komárno', 23),
    @('Košice', 183),
    @('Košice - okolie', 101),
    @('Krupina', 11),
    @('Kysucké Nové Mesto', 4),
    @('Levice', 58),
    @('Levoča', 18),
    @('Liptovský Mikuláš', 29),
    @('Lučenec', 10),
    @('Malacky', 28),
    @('Martin', 1),
    @('Medzilaborce', 5),
    @('Michalovce', 31),
    @('Myjava', 2),
    @('Námestovo', 3),
    @('Nitra', 62),
    @('Nové Mesto nad Váhom', 18),
    @('Nové Zámky', 30),
    @('Partizánske', 13),
    @('Pezinok', 25),
    @('Piešťany', 2),
    @('Poltár', 9),
    @('Poprad', 66),
    @('Považská Bystrica', 29),
    @('Prešov', 8),
    @('Prievidza', 12),
    @('Púchov', 11),
    @('Revúca', 3),
    @('Rimavská Sobota', 38),
    @('Rožňava', 3),
    @('Ružomberok', 35),
    @('Sabinov', 16),
    @('Senec', 46),
    @('Senica', 2),
    @('Skalica', 1),
    @('Snina', 10),
    @('Sobrance', 1),
    @('Spišská Nová Ves', 68),
    @('Stará Ľubovňa', 43),
    @('Stropkov', 14),
    @('Svidník', 25),
    @('Šaľa', 20),
    @('Topoľčany', 33),
    @('Trebišov', 48),
    @('Trenčín', 53),
    @('Trnava', 41),
    @('Turčianske Teplice', 1),
    @('Veľký Krtíš', 9),
    @('Vranov nad Topľou', 36),
    @('Zlaté Moravce', 33),
    @('Zvolen', 21),
    @('Žarnovica', 21),
    @('Žiar nad Hronom', 17),
    @('Žilina', 6)
)

# Insert 71 blank rows above the first row of the existing 44203 block,
# shifting it (and everything below) down.
$ws.Rows("$startRow`:$($startRow + $districts.Length - 1)").Insert()

for ($i = 0; $i -lt $districts.Length; $i++) {
    $row = $startRow + $i
    $pair = $districts[$i]
    $ws.Cells.Item($row, 1).Value2 = $newDate
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($row, 2).Value2 = $pair[0]
    $ws.Cells.Item($row, 3).Value2 = $pair[1]
}
